# Insert a new data row before the current row 76 (shifts existing rows 76-194
# down to 77-195) and populate it with the new observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(76).Insert()

$ws.Range("A76").Value = 1
$ws.Range("B76").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C76").Value = "Arica y Parinacota"
$ws.Range("D76").Value = 44477
$ws.Range("E76").Value = 15
$ws.Range("F76").Value = 100114013
$ws.Range("G76").Value = "Zanahoria"
$ws.Range("H76").Value = "Sin especificar"
$ws.Range("I76").Value = "Primera"
$ws.Range("J76").Value = 80
$ws.Range("K76").Value = 7000
$ws.Range("L76").Value = 8000
$ws.Range("M76").Value = 7500
$ws.Range("N76").Value = "`$/saco 25 kilos"
$ws.Range("O76").Value = "Valle de Camiña"
$ws.Range("P76").Value = 300
$ws.Range("Q76").Value = 25
$ws.Range("R76").Value = "Hortaliza"
